$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-30 Tuesday" "2025-12-31 Wednesday"

Replace-Text "24×83=1992" "22×43=946"
Replace-Text "40×71=2840" "93×30=2790"
Replace-Text "82×73=5986" "20×23=460"
Replace-Text "79×84=6636" "99×17=1683"
Replace-Text "31×47=1457" "52×23=1196"

Replace-Text "13×17=221" "37×41=1517"
Replace-Text "15×18=270" "40×52=2080"
Replace-Text "13×94=1222" "82×48=3936"
Replace-Text "91×70=6370" "92×15=1380"
Replace-Text "65×17=1105" "15×16=240"

Replace-Text "72×61=4392" "49×82=4018"
Replace-Text "67×61=4087" "57×93=5301"
Replace-Text "29×97=2813" "97×89=8633"
Replace-Text "63×60=3780" "54×39=2106"
Replace-Text "17×80=1360" "36×20=720"

Replace-Text "51×73=3723" "69×70=4830"
Replace-Text "91×84=7644" "63×17=1071"
Replace-Text "59×42=2478" "54×73=3942"
Replace-Text "35×17=595" "35×76=2660"
Replace-Text "88×35=3080" "24×88=2112"

Replace-Text "45×53=2385" "64×88=5632"
Replace-Text "91×30=2730" "71×97=6887"
Replace-Text "70×58=4060" "54×18=972"
Replace-Text "39×81=3159" "24×78=1872"
Replace-Text "88×57=5016" "84×81=6804"

Write-Host "Done applying replacements"
